# Weekly update: a new daily price record for Perejil (Feria Lagunitas de
# Puerto Montt) is inserted as the new row 235, pushing the existing rows
# 235:301 down to 236:302 (dimension grows from A1:R301 to A1:R302).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 235; this shifts rows
# 235:301 down to 236:302 and extends the sheet dimension accordingly.
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with the new observation.
$ws.Cells.Item(235, 1).Value = 4
$ws.Cells.Item(235, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235, 3).Value = "Los Lagos"
$ws.Cells.Item(235, 4).Value = 44855
$ws.Cells.Item(235, 5).Value = 10
$ws.Cells.Item(235, 6).Value = 100112044
$ws.Cells.Item(235, 7).Value = "Perejil"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 180
$ws.Cells.Item(235, 11).Value = 5000
$ws.Cells.Item(235, 12).Value = 5000
$ws.Cells.Item(235, 13).Value = 5000
$ws.Cells.Item(235, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(235, 15).Value = "Región Metropolitana"
$ws.Cells.Item(235, 16).Value = 1667
$ws.Cells.Item(235, 17).Value = 3
$ws.Cells.Item(235, 18).Value = "Hortaliza"
